$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp label (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 01:16"

# --- Update country rows: re-ordered countries (re-ranked by totals) + refreshed daily figures ---
$ws.Range("A6").Value = "Estados Unidos"
$ws.Range("B6").Value = 25896
$ws.Range("C6").Value = 6513
$ws.Range("D6").Value = 176
$ws.Range("E6").Value = 25404
$ws.Range("F6").Value = 64
$ws.Range("G6").Value = 60
$ws.Range("H6").Value = 316

$ws.Range("A7").Value = "España"
$ws.Range("B7").Value = 25496
$ws.Range("C7").Value = 3925
$ws.Range("D7").Value = 2125
$ws.Range("E7").Value = 21993
$ws.Range("F7").Value = 1612
$ws.Range("G7").Value = 285
$ws.Range("H7").Value = 1378

$ws.Range("A54").Value = "Panama"
$ws.Range("B54").Value = 245
$ws.Range("C54").Value = 45
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 241
$ws.Range("F54").Value = 7
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 3

$ws.Range("A55").Value = "Sudafrica"
$ws.Range("B55").Value = 240
$ws.Range("C55").Value = 38
$ws.Range("D55").Value = 2
$ws.Range("E55").Value = 238
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0

$ws.Range("A56").Value = "Libano"
$ws.Range("B56").Value = 230
$ws.Range("C56").Value = 53
$ws.Range("D56").Value = 8
$ws.Range("E56").Value = 218
$ws.Range("F56").Value = 4
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 4

$ws.Range("A57").Value = "Irak"
$ws.Range("B57").Value = 214
$ws.Range("C57").Value = 6
$ws.Range("D57").Value = 51
$ws.Range("E57").Value = 146
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 17

$ws.Range("A58").Value = "Croacia"
$ws.Range("B58").Value = 206
$ws.Range("C58").Value = 76
$ws.Range("D58").Value = 5
$ws.Range("E58").Value = 200
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 1

$ws.Range("A59").Value = "Mexico"
$ws.Range("B59").Value = 203
$ws.Range("C59").Value = 39
$ws.Range("D59").Value = 4
$ws.Range("E59").Value = 197
$ws.Range("F59").Value = 1
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 2

$ws.Range("A114").Value = "Paraguay"
$ws.Range("B114").Value = 22
$ws.Range("C114").Value = 4
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 21
$ws.Range("F114").Value = 1
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 1

$ws.Range("A115").Value = "Nigeria"
$ws.Range("B115").Value = 22
$ws.Range("C115").Value = 10
$ws.Range("D115").Value = 1
$ws.Range("E115").Value = 21
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

$ws.Range("A116").Value = "Puerto Rico"
$ws.Range("B116").Value = 21
$ws.Range("C116").Value = 7
$ws.Range("D116").Value = 0
$ws.Range("E116").Value = 21
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

$ws.Range("A118").Value = "Cuba"
$ws.Range("B118").Value = 21
$ws.Range("C118").Value = 0
$ws.Range("D118").Value = 0
$ws.Range("E118").Value = 20
$ws.Range("F118").Value = 0
$ws.Range("G118").Value = 0
$ws.Range("H118").Value = 1

$ws.Range("A119").Value = "Bolivia"
$ws.Range("B119").Value = 19
$ws.Range("C119").Value = 3
$ws.Range("D119").Value = 0
$ws.Range("E119").Value = 19
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 0

$ws.Range("A120").Value = "Jamaica"
$ws.Range("B120").Value = 19
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 2
$ws.Range("E120").Value = 16
$ws.Range("F120").Value = 0
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 1

$ws.Range("A121").Value = "Guayana Francesa"
$ws.Range("B121").Value = 18
$ws.Range("C121").Value = 3
$ws.Range("D121").Value = 0
$ws.Range("E121").Value = 18
$ws.Range("F121").Value = 0
$ws.Range("G121").Value = 0
$ws.Range("H121").Value = 0

$ws.Range("A122").Value = "Monaco"
$ws.Range("B122").Value = 18
$ws.Range("C122").Value = 7
$ws.Range("D122").Value = 1
$ws.Range("E122").Value = 17
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 0
$ws.Range("H122").Value = 0

$ws.Range("A123").Value = "Guyana"
$ws.Range("B123").Value = 18
$ws.Range("C123").Value = 3
$ws.Range("D123").Value = 0
$ws.Range("E123").Value = 17
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0
